$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the species-observation data between row 2 and row 3,
# while leaving shared/common columns (C, D, I, K, P, S, T, U, V, W, Y,
# AA, AD, AE, AG, AT, AW, AX, AY) untouched since those values are
# identical in both rows anyway.

$cols = "A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB"

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
